$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.475.68"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "2.604.58"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.84"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.77"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.07"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.71%  "
$ws.Range("E11").Value = "  -4.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0805"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "3.017.62"
$ws.Range("E15").Value = "  +2.24%  "
$ws.Range("D16").Value = "2.630.15"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "48.486.42"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.93"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.72"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "0.0₃0938"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.35"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.38"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "25.89"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.95"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.64"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  -5.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.06"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.44"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.91"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0791"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.92"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.01"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.25"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.12"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").Value = "2.047.05"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.16"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.85"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.18"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.12"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.78%  "
